# Update report to have fs adjusted weights
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the existing bound tables (IBT, Pension, Retirement) with the new
#    funded-status adjusted weight lower/upper bounds.
# ---------------------------------------------------------------------------

$ibt = $wb.Worksheets.Item("IBT")
$ibt.Range("B6").Value = 0.15
$ibt.Range("C6").Value = 0.55
$ibt.Range("B7").Value = 0.03
$ibt.Range("C7").Value = 0.1
$ibt.Range("B8").Value = 0.03
$ibt.Range("C8").Value = 0.1
$ibt.Range("B9").Value = 0.03
$ibt.Range("C9").Value = 0.1
$ibt.Range("B10").Value = 0.03
$ibt.Range("C10").Value = 0.1
$ibt.Range("B11").Value = 0.02
$ibt.Range("C11").Value = 0.0200000000001

$pension = $wb.Worksheets.Item("Pension")
$pension.Range("B6").Value = 0.25
$pension.Range("C6").Value = 0.55
$pension.Range("B7").Value = 0.05
$pension.Range("C7").Value = 0.12
$pension.Range("B8").Value = 0.05
$pension.Range("C8").Value = 0.12
$pension.Range("B9").Value = 0.03
$pension.Range("C9").Value = 0.1
$pension.Range("B10").Value = 0.04
$pension.Range("C10").Value = 0.12
$pension.Range("B11").Value = 0.02
$pension.Range("C11").Value = 0.0200000000001

$retirement = $wb.Worksheets.Item("Retirement")
$retirement.Range("B6").Value = 0.15
$retirement.Range("C6").Value = 0.35
$retirement.Range("B7").Value = 0.02
$retirement.Range("C7").Value = 0.08
$retirement.Range("B8").Value = 0.02
$retirement.Range("C8").Value = 0.08
$retirement.Range("B9").Value = 0.02
$retirement.Range("C9").Value = 0.08
$retirement.Range("B10").Value = 0.02
$retirement.Range("C10").Value = 0.08
$retirement.Range("B11").Value = 0.02
$retirement.Range("C11").Value = 0.0200000000001

# ---------------------------------------------------------------------------
# 2. Add a new "Unbounded" sheet (after Retirement) carrying the original,
#    wide-open default bounds used before the funded-status adjustment.
# ---------------------------------------------------------------------------

$unbounded = $wb.Worksheets.Add($null, $retirement)
$unbounded.Name = "Unbounded"

$unbounded.Range("A1").Value = "Asset/Liability"
$unbounded.Range("B1").Value = "Lower"
$unbounded.Range("C1").Value = "Upper"

$unbounded.Range("A2").Value = "Liability"
$unbounded.Range("B2").Value = -1.0000000000010001
$unbounded.Range("C2").Value = -0.99999999999999001

$unbounded.Range("A3").Value = "15+ STRIPS"
$unbounded.Range("B3").Value = 0
$unbounded.Range("C3").Value = 1

$unbounded.Range("A4").Value = "Long Corporate"
$unbounded.Range("B4").Value = 0
$unbounded.Range("C4").Value = 1

$unbounded.Range("A5").Value = "Ultra 30Y Futures"
$unbounded.Range("B5").Value = 0
$unbounded.Range("C5").Value = 1

$unbounded.Range("A6").Value = "Equity"
$unbounded.Range("B6").Value = 0
$unbounded.Range("C6").Value = 1

$unbounded.Range("A7").Value = "Liquid Alternatives"
$unbounded.Range("B7").Value = 0
$unbounded.Range("C7").Value = 1

$unbounded.Range("A8").Value = "Private Equity"
$unbounded.Range("B8").Value = 0
$unbounded.Range("C8").Value = 1

$unbounded.Range("A9").Value = "Credit"
$unbounded.Range("B9").Value = 0
$unbounded.Range("C9").Value = 1

$unbounded.Range("A10").Value = "Real Estate"
$unbounded.Range("B10").Value = 0
$unbounded.Range("C10").Value = 1

$unbounded.Range("A11").Value = "Cash"
$unbounded.Range("B11").Value = 0.01
$unbounded.Range("C11").Value = 0.0200000000001

$unbounded.Range("A12").Value = "Hedges"
$unbounded.Range("B12").Value = 0
$unbounded.Range("C12").Value = 0.5

# Match the number formats used on the other bound sheets.
$unbounded.Range("B2").NumberFormat = "0.000000000000"
$unbounded.Range("C2:C12").NumberFormat = "0.00000000000000"

$unbounded.Columns.Item(1).AutoFit()
$unbounded.Columns.Item(2).AutoFit()
$unbounded.Columns.Item(3).AutoFit()

# ---------------------------------------------------------------------------
# 3. Restore the selections / active sheet to match the saved workbook state.
# ---------------------------------------------------------------------------

$unbounded.Range("D9").Select()
$pension.Range("B2:C12").Select()
$retirement.Range("B2:C12").Select()
$ibt.Activate()
$ibt.Range("A9").Select()
